# Update of league bases: swap the odds/result data between each of the
# following adjacent row pairs (the fixtures were recorded swapped -
# columns A, C, D and E -- index, Div, Div Original Name and Date -- stay
# put; column B (id) and the contiguous F..AC block (HomeTeam .. PL_AhUnder)
# trade places between the two rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
  @(13,14),
  @(17,18),
  @(30,31),
  @(46,47),
  @(56,57),
  @(58,59),
  @(62,63),
  @(70,71),
  @(72,73),
  @(74,75),
  @(78,79),
  @(80,81),
  @(82,83),
  @(84,85),
  @(95,96),
  @(106,107),
  @(108,109)
)

# Column B (id) plus the contiguous block F..AC (HomeTeam .. PL_AhUnder).
$cols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

foreach ($pair in $pairs) {
  $r1 = $pair[0]
  $r2 = $pair[1]
  foreach ($c in $cols) {
    $cell1 = $ws.Cells.Item($r1, $c)
    $cell2 = $ws.Cells.Item($r2, $c)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value2 = $v2
    $cell2.Value2 = $v1
  }
}
